$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that carry the row-specific data being rotated between rows 3, 4 and 6
$cols = @("D", "K", "L", "M", "N", "O", "P", "R", "S")

# Snapshot the current ("before") values for the three affected rows
$row3 = @{}
$row4 = @{}
$row6 = @{}
foreach ($col in $cols) {
    $row3[$col] = $ws.Range("$col" + "3").Value()
    $row4[$col] = $ws.Range("$col" + "4").Value()
    $row6[$col] = $ws.Range("$col" + "6").Value()
}

# Apply the cyclic rotation: new row3 = old row6, new row4 = old row3, new row6 = old row4
foreach ($col in $cols) {
    $ws.Range("$col" + "3").Value = $row6[$col]
    $ws.Range("$col" + "4").Value = $row3[$col]
    $ws.Range("$col" + "6").Value = $row4[$col]
}
